# Append new "RightStroke" movement-trial rows for file 00007004_s006_t000
# to the tagged data sheet, mirroring the existing LeftStroke/Normal blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$fileName = "00007004_s006_t000"
$condition = "RightStroke"
$startSamples = 350000
$endSamples = 1700000

$pairs = @(
    @("F3", "F4"),
    @("C3", "C4"),
    @("P3", "P4"),
    @("O1", "O2"),
    @("F7", "F8"),
    @("T3", "T4"),
    @("T5", "T6")
)

$row = 24
foreach ($pair in $pairs) {
    $ws.Range("A10:F10").Copy()
    $ws.Range("A$row`:F$row").PasteSpecial(-4122)

    $ws.Cells.Item(10, 3).Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $fileName
    $ws.Cells.Item($row, 2).Value = $condition
    $ws.Cells.Item($row, 3).Value = $startSamples
    $ws.Cells.Item($row, 4).Value = $endSamples
    $ws.Cells.Item($row, 5).Value = $pair[0]
    $ws.Cells.Item($row, 6).Value = $pair[1]

    $row++
}

$ws.Range("I14").Select()

$wb.Save()
